$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: update product/solution descriptions (Pod Point / Solar / Battery storage text) ---
$f7 = "Home: Pod Point Solo 3S (7kW, Solar integrated).`nOffers: 1. Upfront purchase; 2. ""Plug & Power"" bundle (Hardware discount with Tariff); 3. ""Pod Drive"" (0 upfront, monthly subscription).`nPublic: Pod Point Network (Strategic partner; Tesco/Lidl locations)."
$g7 = "Installer: Contact Solar.`nHardware: Flexible Tier 1 Modular Systems.`nNote: No specific brand advertised (e.g., no Tesla/Powervault). Focus on cost-effective, mid-range solutions."
$h7 = "Installer: Contact Solar (EDF Owned).`nHardware: Standard Tier 1 Monocrystalline (White-label strategy, brands vary by stock)."

$ws.Range("F7").Value2 = $f7
$ws.Range("G7").Value2 = $g7
$ws.Range("H7").Value2 = $h7

# --- Row 7 height increases because of longer wrapped text ---
$ws.Rows.Item(7).RowHeight = 142.5

# --- Row 8: add new F8/G8/H8 content (new hyperlinks + text) ---
$g8 = "Adopts a ""Hardware Agnostic"" strategy: No premium brands (like Tesla Powerwall) are advertised. Focuses on modular, cost-effective solutions installed in-house to maximize margins.`nLink: https://www.contact-solar.co.uk/battery-storage/"
$h8 = "Fulfilled by EDF's subsidiary ""Contact Solar"". Marketing emphasizes technology types (String vs. Micro inverters) rather than specific manufacturers, allowing for flexible Tier 1 hardware sourcing.`nLink: `nhttps://www.contact-solar.co.uk/help-center/our-panels-inverters/"

# F8: plain link cell, display text == URL, no screen tip
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.edfenergy.com/electric-cars/home-charger", [Type]::Missing, [Type]::Missing, "https://www.edfenergy.com/electric-cars/home-charger") | Out-Null

# G8: plain text, no hyperlink
$ws.Range("G8").Value2 = $g8

# H8: hyperlink with long display text + screen tip pointing at the URL
$ws.Hyperlinks.Add($ws.Range("H8"), "https://www.contact-solar.co.uk/help-center/our-panels-inverters/", [Type]::Missing, "https://www.contact-solar.co.uk/help-center/our-panels-inverters/", $h8) | Out-Null

# --- Row 8 height increases to fit new content ---
$ws.Rows.Item(8).RowHeight = 171

# --- Freeze panes: freeze rows 1-3 (header), scroll body to show row 7, active cell H7 ---
$ws.Range("A4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H7").Select()
